$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: extend date header row through column AF (copy format from AC1, then set the value) ---
$ws.Range("AC1").Copy($ws.Range("AD1"))
$ws.Range("AD1").Value = 45318
$ws.Range("AC1").Copy($ws.Range("AE1"))
$ws.Range("AE1").Value = 45319
$ws.Range("AC1").Copy($ws.Range("AF1"))
$ws.Range("AF1").Value = 45320

# --- Row 2 (Pratiksha Bhuse(TL)): AD2/AE2 = Absent, AF2 = Present ---
$ws.Range("AC2").Copy($ws.Range("AD2"))
$ws.Range("AD2").Value = "Absent"
$ws.Range("AC2").Copy($ws.Range("AE2"))
$ws.Range("AE2").Value = "Absent"
$ws.Range("AC2").Copy($ws.Range("AF2"))
$ws.Range("AF2").Value = "Present"

# --- Row 3 (Sangita Survase): AD3/AE3 = Absent, AF3 = Present ---
$ws.Range("AC3").Copy($ws.Range("AD3"))
$ws.Range("AD3").Value = "Absent"
$ws.Range("AC3").Copy($ws.Range("AE3"))
$ws.Range("AE3").Value = "Absent"
$ws.Range("AC3").Copy($ws.Range("AF3"))
$ws.Range("AF3").Value = "Present"

# --- Row 4 (Vaishnavi Wable): AD4/AE4/AF4 = Absent ---
$ws.Range("AC4").Copy($ws.Range("AD4"))
$ws.Range("AD4").Value = "Absent"
$ws.Range("AC4").Copy($ws.Range("AE4"))
$ws.Range("AE4").Value = "Absent"
$ws.Range("AC4").Copy($ws.Range("AF4"))
$ws.Range("AF4").Value = "Absent"

# --- Extend the "Present, Absent,Reason" dropdown validation to cover the new columns ---
# Original sqref was "C5:P5 C2:AC4"; new sqref must become "C5:P5 C2:AF4".
$ws.Range("C5:P5").Validation.Delete()
$ws.Range("C2:AC4").Validation.Delete()
$ws.Range("C2:AF5").Validation.Add(3, 1, 1, '"Present, Absent,Reason"')
$ws.Range("Q5:AF5").Validation.Delete()

# --- Update the active selection to reflect where the user ended up editing ---
$ws.Range("AE10").Select()
